# CranksModelInfo.xlsx edit: replace old rows 9-11 (duplicate Particle-Radius
# timepoints for 2-methylpropanal) and row 16 (the 7.2-radius 2-methylbutanal
# timepoint) with the real fitted-coefficient rows that follow them, shifting
# everything up. Net effect: 4 fewer data rows (25 -> 21), new denominators
# for the per-compound "mass in bean %" / "mass released %" formulas, and a
# refreshed view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 2-methylbutanal block (rows 9-12): was rows 12-15 ------------------
$ws.Range("A9").Value = "2-methylbutanal"
$ws.Range("B9").Value = 13
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 0
$ws.Range("E9").Formula = '=B9/$B$9'
$ws.Range("F9").Formula = '=1-E9'
$ws.Range("G9").Value = 200

$ws.Range("A10").Value = "2-methylbutanal"
$ws.Range("B10").Value = 9.4
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 8
$ws.Range("E10").Formula = '=B10/$B$9'
$ws.Range("F10").Formula = '=1-E10'
$ws.Range("G10").Value = 200

$ws.Range("A11").Value = "2-methylbutanal"
$ws.Range("B11").Value = 8
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 24
$ws.Range("E11").Formula = '=B11/$B$9'
$ws.Range("F11").Formula = '=1-E11'
$ws.Range("G11").Value = 200

$ws.Range("A12").Value = "2-methylbutanal"
$ws.Range("B12").Value = 7.3
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 42
$ws.Range("E12").Formula = '=B12/$B$9'
$ws.Range("F12").Formula = '=1-E12'
$ws.Range("G12").Value = 200

# --- hexanal block (rows 13-21): was rows 17-25, minus the old row 16 ---
$ws.Range("A13").Value = "hexanal"
$ws.Range("B13").Value = 0.7
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 0
$ws.Range("E13").Formula = '=B13/$B$13'
$ws.Range("F13").Formula = '=1-E13'
$ws.Range("G13").Value = 200

$ws.Range("A14").Value = "hexanal"
$ws.Range("B14").Value = 0.45
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 3
$ws.Range("E14").Formula = '=B14/$B$13'
$ws.Range("F14").Formula = '=1-E14'
$ws.Range("G14").Value = 200

$ws.Range("A15").Value = "hexanal"
$ws.Range("B15").Value = 0.52
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 5
$ws.Range("E15").Formula = '=B15/$B$13'
$ws.Range("F15").Formula = '=1-E15'
$ws.Range("G15").Value = 200

$ws.Range("A16").Value = "hexanal"
$ws.Range("B16").Value = 0.43
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 8
$ws.Range("E16").Formula = '=B16/$B$13'
$ws.Range("F16").Formula = '=1-E16'
$ws.Range("G16").Value = 200

$ws.Range("A17").Value = "hexanal"
$ws.Range("B17").Value = 0.42
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 12
$ws.Range("E17").Formula = '=B17/$B$13'
$ws.Range("F17").Formula = '=1-E17'
$ws.Range("G17").Value = 200

$ws.Range("A18").Value = "hexanal"
$ws.Range("B18").Value = 0.41
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 16
$ws.Range("E18").Formula = '=B18/$B$13'
$ws.Range("F18").Formula = '=1-E18'
$ws.Range("G18").Value = 200

$ws.Range("A19").Value = "hexanal"
$ws.Range("B19").Value = 0.39
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 24
$ws.Range("E19").Formula = '=B19/$B$13'
$ws.Range("F19").Formula = '=1-E19'
$ws.Range("G19").Value = 200

$ws.Range("A20").Value = "hexanal"
$ws.Range("B20").Value = 0.36
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 28
$ws.Range("E20").Formula = '=B20/$B$13'
$ws.Range("F20").Formula = '=1-E20'
$ws.Range("G20").Value = 200

$ws.Range("A21").Value = "hexanal"
$ws.Range("B21").Value = 0.34
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 40
$ws.Range("E21").Formula = '=B21/$B$13'
$ws.Range("F21").Formula = '=1-E21'
$ws.Range("G21").Value = 200

# --- rows 22-25 no longer exist; clear them out (shrinks dimension/used
#     range back down to G21 without disturbing the now-correct rows above)
$ws.Range("A22:G25").ClearContents()

# --- view: scroll back to the top-left (drops the stale topLeftCell="A4")
#     and leave the selection where the author's last save left it
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E25").Select()
